$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.647.61'
$ws.Range("D3").Value = '1.598.16'
$ws.Range("E3").Value = '  +0.14%  '
$ws.Range("E4").Value = '  +0.15%  '
$ws.Range("D5").Value = "'" + '211.32'
$ws.Range("E5").Value = '  +0.11%  '
$ws.Range("E6").Value = '  +0.77%  '
$ws.Range("E7").Value = '  +0.12%  '
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("E9").Value = '  +0.10%  '
$ws.Range("D10").Value = "'" + '19.53'
$ws.Range("E10").Value = '  -0.65%  '
$ws.Range("D11").Value = "'" + '0.0839'
$ws.Range("E11").Value = '  +0.36%  '
$ws.Range("D12").Value = '1.822.28'
$ws.Range("E12").Value = '  +0.16%  '
$ws.Range("D13").Value = '1.615.61'
$ws.Range("E13").Value = '  +1.39%  '
$ws.Range("E14").Value = '  +0.08%  '
$ws.Range("D15").Value = "'" + '0.524'
$ws.Range("E15").Value = '  +0.22%  '
$ws.Range("D17").Value = '26.635.85'
$ws.Range("D18").Value = '0.0₃0734'
$ws.Range("E18").Value = '  +0.61%  '
$ws.Range("E19").Value = '  +0.10%  '
$ws.Range("D20").Value = "'" + '207.94'
$ws.Range("E20").Value = '  -0.58%  '
$ws.Range("E21").Value = '  +5.55%  '
$ws.Range("E22").Value = '  +0.76%  '
$ws.Range("E23").Value = '  +0.88%  '
$ws.Range("E24").Value = '  +0.34%  '
$ws.Range("D25").Value = "'" + '145.46'
$ws.Range("E25").Value = '  -0.80%  '
$ws.Range("E26").Value = '  +0.09%  '
$ws.Range("E27").Value = '  +0.03%  '
$ws.Range("E28").Value = '  +0.24%  '
$ws.Range("E29").Value = '  -0.01%  '
$ws.Range("D30").Value = "'" + '0.0512'
$ws.Range("E30").Value = '  +2.01%  '
$ws.Range("E32").Value = '  +0.37%  '
$ws.Range("E33").Value = '  +0.96%  '
$ws.Range("B34").Value = 'ImmutableX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D34").Value = "'" + '0.623'
$ws.Range("E34").Value = '  -9.09%  '
$ws.Range("B35").Value = 'Maker'
$ws.Range("C35").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D35").Value = '1.274.72'
$ws.Range("E35").Value = '  -1.42%  '
$ws.Range("E36").Value = '  +0.70%  '
$ws.Range("E37").Value = '  +0.60%  '
$ws.Range("D38").Value = "'" + '0.0171'
$ws.Range("E38").Value = '  -0.42%  '
$ws.Range("E39").Value = '  +20.22%  '
$ws.Range("E40").Value = '  -0.08%  '
$ws.Range("E41").Value = '  +2.87%  '
$ws.Range("E42").Value = '  +0.40%  '
$ws.Range("D43").Value = "'" + '0.784'
$ws.Range("E43").Value = '  -0.93%  '
$ws.Range("D44").Value = "'" + '63.91'
$ws.Range("E44").Value = '  +0.57%  '
$ws.Range("D45").Value = '1.734.91'
$ws.Range("E45").Value = '  +0.25%  '
$ws.Range("D46").Value = "'" + '90.16'
$ws.Range("E46").Value = '  +0.50%  '
$ws.Range("E47").Value = '  -2.49%  '
$ws.Range("E48").Value = '  +3.29%  '
$ws.Range("E49").Value = '  +0.96%  '
$ws.Range("D51").Value = "'" + '7.39'
$ws.Range("E51").Value = '  -1.28%  '
